$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# e007 "Morning Briefing - Weather Roll" entry (row 8, column B):
# append an extra trailing <LineBreak/> to the rich-text description.
$newB8 = @'
<Bold>e007 Morning Briefing - Weather Roll</Bold> <InlineUIContainer><Button Content='r4.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The 
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
 Table determines weather for today:  
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@

# e008 "Type of Snow" entry (row 9, column B):
# add a space before the line break after "on the", and append an extra
# trailing <LineBreak/>.
$newB9 = @'
<Bold>e008 Type of Snow</Bold> 
<LineBreak/><LineBreak/>
Snow is in the forecast. Roll for type of snow on the 
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
 Table:  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@

# Here-strings in PowerShell keep the trailing newline before the closing
# '@ marker; trim it so the text matches the original (no trailing \n).
$newB8 = $newB8.TrimEnd("`r", "`n")
$newB9 = $newB9.TrimEnd("`r", "`n")

$ws.Range("B8").Value = $newB8
$ws.Range("B9").Value = $newB9

# Move the selection/scroll position to reflect where the edit was made.
$ws.Range("B8").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
